$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The row containing the lone "RACE" header label (row 11) is deleted.
# Deleting the entire row shifts subsequent rows up, and Excel will
# automatically drop the now-unused "RACE" shared string on save.
$ws.Rows.Item(11).Delete()

# Update the active selection/cell to match the post-edit state.
$ws.Range("A11:XFD11").Select()
